# Updates the cryptos list (Coin/Link/Price/Volume(1h)) to match the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Columns B (Coin) and C (Link) are plain text, so they are safe to set
# directly via .Value. Columns D (Price) and E (Volume(1h)) are stored as
# text in the workbook even though many of the values look numeric (e.g.
# "35.855.46", "0.623", "1.01") -- Excel would otherwise auto-convert a
# plain .Value assignment like that into a real number and mangle the
# display (dropping the thousands dots, trailing zeros, etc). To avoid
# that, each D/E write temporarily forces the cell to Text format ("@"),
# assigns the literal string, then restores the cells original Style so
# no formatting/style footprint is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$__style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.855.46"
$ws.Range("D2").Style = $__style
$__style = $ws.Range("E2").Style
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -4.70%  "
$ws.Range("E2").Style = $__style

# Row 3
$__style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.952.72"
$ws.Range("D3").Style = $__style
$__style = $ws.Range("E3").Style
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.39%  "
$ws.Range("E3").Style = $__style

# Row 4
$__style = $ws.Range("E4").Style
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("E4").Style = $__style

# Row 5
$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.53"
$ws.Range("D5").Style = $__style
$__style = $ws.Range("E5").Style
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -4.52%  "
$ws.Range("E5").Style = $__style

# Row 6
$__style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.623"
$ws.Range("D6").Style = $__style
$__style = $ws.Range("E6").Style
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.68%  "
$ws.Range("E6").Style = $__style

# Row 7
$__style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.17"
$ws.Range("D7").Style = $__style
$__style = $ws.Range("E7").Style
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -9.42%  "
$ws.Range("E7").Style = $__style

# Row 8
$__style = $ws.Range("E8").Style
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("E8").Style = $__style

# Row 9
$__style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.368"
$ws.Range("D9").Style = $__style
$__style = $ws.Range("E9").Style
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.25%  "
$ws.Range("E9").Style = $__style

# Row 10
$__style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "56.15"
$ws.Range("D10").Style = $__style
$__style = $ws.Range("E10").Style
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.57%  "
$ws.Range("E10").Style = $__style

# Row 11
$__style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0793"
$ws.Range("D11").Style = $__style
$__style = $ws.Range("E11").Style
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.18%  "
$ws.Range("E11").Style = $__style

# Row 12
$__style = $ws.Range("E12").Style
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("E12").Style = $__style

# Row 13
$__style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.842"
$ws.Range("D13").Style = $__style
$__style = $ws.Range("E13").Style
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -7.06%  "
$ws.Range("E13").Style = $__style

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$__style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "13.81"
$ws.Range("D14").Style = $__style
$__style = $ws.Range("E14").Style
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -8.55%  "
$ws.Range("E14").Style = $__style

# Row 15
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$__style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.42"
$ws.Range("D15").Style = $__style
$__style = $ws.Range("E15").Style
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.78%  "
$ws.Range("E15").Style = $__style

# Row 16
$__style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.247.32"
$ws.Range("D16").Style = $__style
$__style = $ws.Range("E16").Style
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.02%  "
$ws.Range("E16").Style = $__style

# Row 17
$__style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.34"
$ws.Range("D17").Style = $__style
$__style = $ws.Range("E17").Style
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.93%  "
$ws.Range("E17").Style = $__style

# Row 18
$__style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.961.04"
$ws.Range("D18").Style = $__style
$__style = $ws.Range("E18").Style
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.69%  "
$ws.Range("E18").Style = $__style

# Row 19
$__style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "35.804.49"
$ws.Range("D19").Style = $__style
$__style = $ws.Range("E19").Style
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.38%  "
$ws.Range("E19").Style = $__style

# Row 20
$__style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.49"
$ws.Range("D20").Style = $__style
$__style = $ws.Range("E20").Style
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.90%  "
$ws.Range("E20").Style = $__style

# Row 21
$__style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0846"
$ws.Range("D21").Style = $__style
$__style = $ws.Range("E21").Style
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.15%  "
$ws.Range("E21").Style = $__style

# Row 22
$__style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "234.75"
$ws.Range("D22").Style = $__style
$__style = $ws.Range("E22").Style
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("E22").Style = $__style

# Row 23
$__style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.14"
$ws.Range("D23").Style = $__style
$__style = $ws.Range("E23").Style
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.90%  "
$ws.Range("E23").Style = $__style

# Row 24
$__style = $ws.Range("E24").Style
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("E24").Style = $__style

# Row 25
$__style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.52"
$ws.Range("D25").Style = $__style
$__style = $ws.Range("E25").Style
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.99%  "
$ws.Range("E25").Style = $__style

# Row 26
$__style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.27"
$ws.Range("D26").Style = $__style
$__style = $ws.Range("E26").Style
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.24%  "
$ws.Range("E26").Style = $__style

# Row 27
$__style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.64"
$ws.Range("D27").Style = $__style
$__style = $ws.Range("E27").Style
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("E27").Style = $__style

# Row 28
$__style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "158.85"
$ws.Range("D28").Style = $__style
$__style = $ws.Range("E28").Style
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.94%  "
$ws.Range("E28").Style = $__style

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$__style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.63"
$ws.Range("D29").Style = $__style
$__style = $ws.Range("E29").Style
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("E29").Style = $__style

# Row 30
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$__style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.129"
$ws.Range("D30").Style = $__style
$__style = $ws.Range("E30").Style
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +16.32%  "
$ws.Range("E30").Style = $__style

# Row 31
$__style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.118"
$ws.Range("D31").Style = $__style
$__style = $ws.Range("E31").Style
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.82%  "
$ws.Range("E31").Style = $__style

# Row 32
$__style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.12"
$ws.Range("D32").Style = $__style
$__style = $ws.Range("E32").Style
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -7.16%  "
$ws.Range("E32").Style = $__style

# Row 33
$__style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.78"
$ws.Range("D33").Style = $__style
$__style = $ws.Range("E33").Style
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -8.51%  "
$ws.Range("E33").Style = $__style

# Row 34
$__style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0612"
$ws.Range("D34").Style = $__style
$__style = $ws.Range("E34").Style
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("E34").Style = $__style

# Row 35
$__style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.34"
$ws.Range("D35").Style = $__style
$__style = $ws.Range("E35").Style
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -8.44%  "
$ws.Range("E35").Style = $__style

# Row 36
$__style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.23"
$ws.Range("D36").Style = $__style
$__style = $ws.Range("E36").Style
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.16%  "
$ws.Range("E36").Style = $__style

# Row 37
$__style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.01"
$ws.Range("D37").Style = $__style
$__style = $ws.Range("E37").Style
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("E37").Style = $__style

# Row 38
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$__style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.83"
$ws.Range("D38").Style = $__style
$__style = $ws.Range("E38").Style
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("E38").Style = $__style

# Row 39
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$__style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.24"
$ws.Range("D39").Style = $__style
$__style = $ws.Range("E39").Style
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -8.57%  "
$ws.Range("E39").Style = $__style

# Row 40
$__style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.00"
$ws.Range("D40").Style = $__style
$__style = $ws.Range("E40").Style
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +10.16%  "
$ws.Range("E40").Style = $__style

# Row 41
$__style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0977"
$ws.Range("D41").Style = $__style
$__style = $ws.Range("E41").Style
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.67%  "
$ws.Range("E41").Style = $__style

# Row 42
$__style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.21"
$ws.Range("D42").Style = $__style
$__style = $ws.Range("E42").Style
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.41%  "
$ws.Range("E42").Style = $__style

# Row 43
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$__style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.83"
$ws.Range("D43").Style = $__style
$__style = $ws.Range("E43").Style
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.73%  "
$ws.Range("E43").Style = $__style

# Row 44
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$__style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0209"
$ws.Range("D44").Style = $__style
$__style = $ws.Range("E44").Style
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.26%  "
$ws.Range("E44").Style = $__style

# Row 45
$__style = $ws.Range("E45").Style
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.38%  "
$ws.Range("E45").Style = $__style

# Row 46
$__style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.29"
$ws.Range("D46").Style = $__style
$__style = $ws.Range("E46").Style
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.12%  "
$ws.Range("E46").Style = $__style

# Row 47
$__style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.85"
$ws.Range("D47").Style = $__style
$__style = $ws.Range("E47").Style
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.48%  "
$ws.Range("E47").Style = $__style

# Row 48
$__style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.45"
$ws.Range("D48").Style = $__style
$__style = $ws.Range("E48").Style
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -7.71%  "
$ws.Range("E48").Style = $__style

# Row 49
$__style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.328.48"
$ws.Range("D49").Style = $__style
$__style = $ws.Range("E49").Style
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -6.86%  "
$ws.Range("E49").Style = $__style

# Row 50
$__style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.78"
$ws.Range("D50").Style = $__style
$__style = $ws.Range("E50").Style
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.32%  "
$ws.Range("E50").Style = $__style

# Row 51
$__style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.140.74"
$ws.Range("D51").Style = $__style
$__style = $ws.Range("E51").Style
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.01%  "
$ws.Range("E51").Style = $__style

